# Add code for experiment 3, MASS for T4 using T1, T2, T3
# The sheet now gets a title label in A1 (the sheet/experiment name),
# and the active selection is moved to A9 (ready for the next block of
# results to be typed below the existing table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Label the results table with the experiment/sheet name.
$ws.Range("A1").Value = "AMC_5"

# Leave the selection where the next experiment's data will be entered.
$ws.Range("A9").Select() | Out-Null
